# FanPowerAllowances-T24N_2022.xlsx -- add new Exhaust Systems Base Allowance /
# Supply+Return ERV / Return Filter columns (J:M) to the first table, shifting
# the old SZVAV column out to N, and refresh the "Single Zone VAV..." column
# note that used to live (mis-filed) under the now-removed "Energy Recovery"
# string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: top header band (wrapped descriptive text, style carries over automatically) ---
$ws.Range("J3").Value = "Exhaust Systems Base Allowance"
$ws.Range("K3").Value = "Supply - Energy Recovery (Enthalpy Recovery Ratio ? 0.60 and <0.65)  "
$ws.Range("L3").Value = "Return - Energy Recovery (Enthalpy Recovery Ratio ? 0.60 and <0.65)  "
$ws.Range("M3").Value = "Return - Filter (any MERV value)"
$ws.Range("N3").Value = "Single Zone VAV Systems that are capable of turning down to 50% of full load airflow at a maximum of 30% design wattage"

# --- Row 4: variable-name header band ---
$ws.Range("J4").Value = "ExhBaseAllow"
$ws.Range("K4").Value = "ERVSupply"
$ws.Range("L4").Value = "ERVReturn"
$ws.Range("M4").Value = "RetFilter"
$ws.Range("N4").Value = "SZVAV"
$ws.Range("J4:M4").WrapText = $true

# --- Data rows 5-10: new ExhBaseAllow/ERVSupply/ERVReturn/RetFilter values,
#     old SZVAV data (previously in column K) now lives in column N ---
$ws.Range("J5").Value = 0.221
$ws.Range("K5").Value = 0.184
$ws.Range("L5").Value = 0.19
$ws.Range("M5").Value = 0.046
$ws.Range("N5").Value = 0
$ws.Range("N5").NumberFormat = "0.000"

$ws.Range("J6").Value = 0.246
$ws.Range("K6").Value = 0.155
$ws.Range("L6").Value = 0.163
$ws.Range("M6").Value = 0.041
$ws.Range("N6").Value = 0
$ws.Range("N6").NumberFormat = "0.000"

$ws.Range("J7").Value = 0.236
$ws.Range("K7").Value = 0.144
$ws.Range("L7").Value = 0.146
$ws.Range("M7").Value = 0.036
$ws.Range("N7").Value = 0
$ws.Range("N7").NumberFormat = "0.000"

$ws.Range("J8").Value = 0.186
$ws.Range("K8").Value = 0.19
$ws.Range("L8").Value = 0.191
$ws.Range("M8").Value = 0.046
$ws.Range("N8").Value = 0.07
$ws.Range("N8").NumberFormat = "0.000"

$ws.Range("J9").Value = 0.184
$ws.Range("K9").Value = 0.163
$ws.Range("L9").Value = 0.166
$ws.Range("M9").Value = 0.041
$ws.Range("N9").Value = 0.1
$ws.Range("N9").NumberFormat = "0.000"

$ws.Range("J10").Value = 0.19
$ws.Range("K10").Value = 0.146
$ws.Range("L10").Value = 0.148
$ws.Range("M10").Value = 0.036
$ws.Range("N10").Value = 0.089
$ws.Range("N10").NumberFormat = "0.000"

# --- Column widths: J:N get the wider 22.71-char width used by the new columns ---
$ws.Range("J1:N1").ColumnWidth = 22

# --- Row heights for the widened header rows ---
$ws.Rows.Item(3).RowHeight = 90
$ws.Rows.Item(4).RowHeight = 30

# --- Selection moves to M15 (cursor location when file was last saved) ---
$ws.Range("M15").Select()
